$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Skill Description") before the existing SFIA Level column.
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "Skill Description"

# Map SkillCode -> human readable skill description
$skillDescriptions = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "SORC"       = "Sourcing"
    "SUPP"       = "Supplier management"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $skillDescriptions[$code]
}
